# Add the two new "Diapositive6b" / "Diapositive6c" instruction rows that
# Rita supplied, right after the existing "Diapositive6" row (old row 7)
# and before the "Alors, on commence !" block (old row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room: push the old rows 8-16 down to 10-18.
$ws.Rows("8:9").Insert()

# Write the two new "Img" names first so they land in the shared-string
# table before the two long instruction texts (matches the original
# authoring order: Diapositive6b, Diapositive6c, then the two texts).
$ws.Range("C8").Value = "Diapositive6b"
$ws.Range("C9").Value = "Diapositive6c"

# New row 8 - Diapositive6b
$ws.Range("A8").Value = "Instruct"
$ws.Range("B8").Value = "Dans ce qui suit tu vas bien entendre les sauts`net essayer d'appuyer sur la bonne touche."
$ws.Range("B8").WrapText = $true
$ws.Range("D8").Value = "Key"
$ws.Range("E8").Value = "None"
$ws.Rows(8).RowHeight = 31.5

# New row 9 - Diapositive6c
$ws.Range("A9").Value = "Instruct"
$ws.Range("B9").Value = "Si la première fois tu n'as pas bien entendu les sauts, tu peux réécouter`nen appuyant sur ESPACE"
$ws.Range("B9").WrapText = $true
$ws.Range("D9").Value = "Key"
$ws.Range("E9").Value = "None"
$ws.Rows(9).RowHeight = 31.5

# The rows that used to be 8-16 kept their old (now stale) explicit row
# heights after the insert; re-fit them so they fall back to the sheet's
# default height, same as Excel leaves them once re-saved.
$ws.Rows("10:14").AutoFit()

# Select the two freshly-added rows, same as the author left the workbook.
[void]$ws.Range("A8:E9").Select()
